$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("110:110").Delete()
